$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Range("G1").Value = "Percentage Change 2019-20"

# Data values for the new column (percentage change 2019-20)
$ws.Range("G2").Value = -2.8
$ws.Range("G3").Value = -3.5
$ws.Range("G4").Value = -5.9
$ws.Range("G5").Value = -2.9
$ws.Range("G6").Value = -4.0
$ws.Range("G7").Value = -4.4
$ws.Range("G8").Value = -4.0
$ws.Range("G9").Value = -4.1
$ws.Range("G10").Value = -2.5
$ws.Range("G11").Value = -0.7
$ws.Range("G12").Value = -3.8
$ws.Range("G13").Value = -2.5
$ws.Range("G14").Value = -2.5
$ws.Range("G15").Value = -5.4
$ws.Range("G16").Value = -5.4
$ws.Range("G17").Value = -1.5
$ws.Range("G18").Value = -3.7
$ws.Range("G19").Value = -4.9
$ws.Range("G20").Value = -3.1
$ws.Range("G21").Value = -0.9
$ws.Range("G22").Value = -4.5
$ws.Range("G23").Value = -5.4
$ws.Range("G24").Value = -4.1
$ws.Range("G25").Value = -5.5
$ws.Range("G26").Value = -2.8
$ws.Range("G27").Value = -4.1
$ws.Range("G28").Value = -2.7
$ws.Range("G29").Value = -3.7
$ws.Range("G30").Value = -6.1
$ws.Range("G31").Value = -2.3
$ws.Range("G32").Value = -0.1
$ws.Range("G33").Value = -3.0
$ws.Range("G34").Value = -4.6
$ws.Range("G35").Value = -4.9
$ws.Range("G36").Value = -2.1
$ws.Range("G37").Value = -5.4
$ws.Range("G38").Value = -3.1
$ws.Range("G39").Value = -8.0
$ws.Range("G40").Value = -4.7
$ws.Range("G41").Value = -5.5
$ws.Range("G42").Value = -1.1
$ws.Range("G43").Value = -3.9
$ws.Range("G44").Value = -3.7
$ws.Range("G45").Value = -4.5
$ws.Range("G46").Value = -3.5
$ws.Range("G47").Value = -2.7
$ws.Range("G48").Value = -1.7
$ws.Range("G49").Value = -3.0
$ws.Range("G50").Value = -7.0
$ws.Range("G51").Value = -5.4

# Apply font size 14 to the new column (header + data)
$ws.Range("G1:G51").Font.Size = 14

# Set the column width to match the source
$ws.Columns.Item(7).ColumnWidth = 32.166666666666664

# Select the whole new column, matching the final selection state
$ws.Columns.Item(7).Select()
